$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 32
$ws1.Range("F3").Value = 21109
$ws1.Range("F6").Value = 1125
$ws1.Range("F8").Value = 7878
$ws1.Range("F10").Value = 39
$ws1.Range("F11").Value = 761
$ws1.Range("F12").Value = 304
$ws1.Range("F15").Value = 162
$ws1.Range("F18").Value = 223
$ws1.Range("F20").Value = 512
$ws1.Range("F21").Value = 81
$ws1.Range("F22").Value = 706
$ws1.Range("F26").Value = 345
$ws1.Range("F27").Value = 1174
$ws1.Range("F30").Value = 217
$ws1.Range("F33").Value = 3
$ws1.Range("F34").Value = 131
$ws1.Range("F35").Value = 5015
$ws1.Range("F36").Value = 31
$ws1.Range("F38").Value = 39
$ws1.Range("F40").Value = 13033
$ws1.Range("F42").Value = 125
$ws1.Range("F45").Value = 302
$ws1.Range("F46").Value = 420
$ws1.Range("F47").Value = 4055
$ws1.Range("F48").Value = 329

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 32
$ws4.Range("F3").Value = 21109
$ws4.Range("F5").Value = 1125
$ws4.Range("F7").Value = 7878
$ws4.Range("F9").Value = 39
$ws4.Range("F10").Value = 761
$ws4.Range("F11").Value = 304
$ws4.Range("F14").Value = 162
$ws4.Range("F16").Value = 223
$ws4.Range("F18").Value = 512
$ws4.Range("F19").Value = 81
$ws4.Range("F20").Value = 706
$ws4.Range("F24").Value = 345
$ws4.Range("F25").Value = 1174
$ws4.Range("F28").Value = 217
$ws4.Range("F32").Value = 3
$ws4.Range("F33").Value = 131
$ws4.Range("F35").Value = 5015
$ws4.Range("F36").Value = 31
$ws4.Range("F38").Value = 39
$ws4.Range("F40").Value = 13033
$ws4.Range("F42").Value = 125
$ws4.Range("F45").Value = 302
$ws4.Range("F46").Value = 420
$ws4.Range("F47").Value = 4055
$ws4.Range("F48").Value = 329
